# Updated cryptos list on Tue Apr 18 11:23:21 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "30.112.65"
$ws.Cells.Item(2, 5).Value = "  +0.14%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.117.33"
$ws.Cells.Item(3, 5).Value = "  +0.63%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.19%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "346.40"
$ws.Cells.Item(5, 5).Value = "  +0.39%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.14%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.5209"
$ws.Cells.Item(7, 5).Value = "  +0.54%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.4470"
$ws.Cells.Item(8, 5).Value = "  -0.14%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "54.29"
$ws.Cells.Item(9, 5).Value = "  +3.81%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.09364"
$ws.Cells.Item(10, 5).Value = "  -1.36%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.185"
$ws.Cells.Item(11, 5).Value = "  +0.69%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "25.32"
$ws.Cells.Item(12, 5).Value = "  +0.36%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "8.683"
$ws.Cells.Item(13, 5).Value = "  +7.44%  "

# Row 14
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value = "2.136.61"
$ws.Cells.Item(14, 5).Value = "  +1.26%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "Polkadot"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "6.974"
$ws.Cells.Item(15, 5).Value = "  +3.41%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "102.61"
$ws.Cells.Item(16, 5).Value = "  +3.13%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.03%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "1.007"
$ws.Cells.Item(18, 5).Value = "  -0.23%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +4.58%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.06702"
$ws.Cells.Item(20, 5).Value = "  -0.06%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.309"
$ws.Cells.Item(21, 5).Value = "  +1.96%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.006"
$ws.Cells.Item(22, 5).Value = "  -0.09%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "30.136.92"
$ws.Cells.Item(23, 5).Value = "  -0.04%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "12.75"
$ws.Cells.Item(24, 5).Value = "  +0.45%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.330"
$ws.Cells.Item(25, 5).Value = "  +0.61%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "2.356.10"
$ws.Cells.Item(26, 5).Value = "  +0.00%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +0.66%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.552"
$ws.Cells.Item(28, 5).Value = "  +0.78%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "163.18"
$ws.Cells.Item(29, 5).Value = "  -0.79%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "134.27"
$ws.Cells.Item(30, 5).Value = "  +0.54%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.157"
$ws.Cells.Item(31, 5).Value = "  -0.15%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.787"
$ws.Cells.Item(32, 5).Value = "  +9.87%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +0.18%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.822"
$ws.Cells.Item(34, 5).Value = "  +10.61%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Filecoin"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "6.289"
$ws.Cells.Item(35, 5).Value = "  +0.59%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.966"
$ws.Cells.Item(36, 5).Value = "  +0.44%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "10.76"
$ws.Cells.Item(37, 5).Value = "  +6.19%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +2.86%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.06880"
$ws.Cells.Item(39, 5).Value = "  +1.44%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.7128"
$ws.Cells.Item(40, 5).Value = "  +2.55%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "12.76"
$ws.Cells.Item(41, 5).Value = "  +2.37%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.2250"
$ws.Cells.Item(42, 5).Value = "  -1.33%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +1.85%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.6974"
$ws.Cells.Item(44, 5).Value = "  +3.91%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "14.70"
$ws.Cells.Item(45, 5).Value = "  +2.96%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.402"
$ws.Cells.Item(46, 5).Value = "  +5.29%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.006"
$ws.Cells.Item(47, 5).Value = "  +0.24%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.265"
$ws.Cells.Item(48, 5).Value = "  +7.75%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.624"
$ws.Cells.Item(49, 5).Value = "  -0.45%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.00000000351"
$ws.Cells.Item(50, 5).Value = "  -0.58%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "ThetaToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.219"
$ws.Cells.Item(51, 5).Value = "  +9.35%  "
